$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated case-study results (B2:D5)
$ws.Range("B2").Value = 89289.091942414307
$ws.Range("C2").Value = 61900.571533217502
$ws.Range("D2").Value = 205.22809740707601

$ws.Range("B3").Value = 86910.332389806397
$ws.Range("C3").Value = 59554.102570634597
$ws.Range("D3").Value = 237.51868743218799

$ws.Range("B4").Value = 85104.733203326701
$ws.Range("C4").Value = 57788.213950614001
$ws.Range("D4").Value = 277.22925389091398

$ws.Range("B5").Value = 83714.160106570795
$ws.Range("C5").Value = 56459.615119203401
$ws.Range("D5").Value = 339.20351923521901

# Move the active selection from G9 to F11
$ws.Range("F11").Select()
